$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2: update D2, E2
Set-TextValue "D2" "26.296.00"
Set-TextValue "E2" "  +0.54%  "

# Row 3: update D3, E3
Set-TextValue "D3" "1.594.86"
Set-TextValue "E3" "  +0.98%  "

# Row 4: update D4, E4
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.34%  "

# Row 5: update D5, E5
Set-TextValue "D5" "212.81"
Set-TextValue "E5" "  +1.67%  "

# Row 7: update E7
Set-TextValue "E7" "  -0.31%  "

# Row 8: update E8
Set-TextValue "E8" "  +0.41%  "

# Row 9: update D9, E9
Set-TextValue "D9" "0.0609"
Set-TextValue "E9" "  +0.02%  "

# Row 10: update D10, E10
Set-TextValue "D10" "19.39"
Set-TextValue "E10" "  -0.68%  "

# Row 11: update D11, E11
Set-TextValue "D11" "0.0849"
Set-TextValue "E11" "  +0.54%  "

# Row 12: update D12, E12
Set-TextValue "D12" "1.815.52"
Set-TextValue "E12" "  +0.79%  "

# Row 13: update D13, E13
Set-TextValue "D13" "1.577.41"
Set-TextValue "E13" "  -0.55%  "

# Row 14: update E14
Set-TextValue "E14" "  -0.05%  "

# Row 15: update D15, E15
Set-TextValue "D15" "0.522"
Set-TextValue "E15" "  +1.54%  "

# Row 16: update D16, E16
Set-TextValue "D16" "64.47"
Set-TextValue "E16" "  +0.02%  "

# Row 17: update D17, E17
Set-TextValue "D17" "26.282.20"
Set-TextValue "E17" "  +0.47%  "

# Row 18: update D18, E18
Set-TextValue "D18" "0.0₃0728"
Set-TextValue "E18" "  -0.53%  "

# Row 19: update D19, E19
Set-TextValue "D19" "7.45"
Set-TextValue "E19" "  +2.57%  "

# Row 20: update D20, E20
Set-TextValue "D20" "213.81"
Set-TextValue "E20" "  +2.95%  "

# Row 21: update E21
Set-TextValue "E21" "  -0.26%  "

# Row 22: update D22, E22
Set-TextValue "D22" "4.29"
Set-TextValue "E22" "  +0.93%  "

# Row 23: update D23, E23
Set-TextValue "D23" "9.03"
Set-TextValue "E23" "  +1.71%  "

# Row 24: update E24
Set-TextValue "E24" "  -2.55%  "

# Row 25: update D25, E25
Set-TextValue "D25" "144.65"
Set-TextValue "E25" "  +0.16%  "

# Row 26: update E26
Set-TextValue "E26" "  -0.39%  "

# Row 27: update D27, E27
Set-TextValue "D27" "7.06"
Set-TextValue "E27" "  +1.06%  "

# Row 28: update E28
Set-TextValue "E28" "  -0.52%  "

# Row 29: update D29, E29
Set-TextValue "D29" "15.22"
Set-TextValue "E29" "  +0.10%  "

# Row 30: update D30, E30
Set-TextValue "D30" "0.0499"
Set-TextValue "E30" "  -1.25%  "

# Row 31: update E31
Set-TextValue "E31" "  +1.06%  "

# Row 32: update E32
Set-TextValue "E32" "  -0.18%  "

# Row 33: update E33
Set-TextValue "E33" "  -0.12%  "

# Row 34: update D34, E34
Set-TextValue "D34" "1.339.68"
Set-TextValue "E34" "  +4.99%  "

# Row 35: update E35
Set-TextValue "E35" "  -1.00%  "

# Row 36: update E36
Set-TextValue "E36" "  -0.27%  "

# Row 37: update D37, E37
Set-TextValue "D37" "0.591"
Set-TextValue "E37" "  -3.24%  "

# Row 38: update E38
Set-TextValue "E38" "  +0.29%  "

# Row 39: update D39, E39
Set-TextValue "D39" "0.822"
Set-TextValue "E39" "  +0.63%  "

# Row 40: update D40, E40
Set-TextValue "D40" "1.02"
Set-TextValue "E40" "  -11.58%  "

# Row 41: update B41, C41, D41, E41
Set-TextValue "B41" "FraxShare"
Set-TextValue "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "5.72"
Set-TextValue "E41" "  +2.93%  "

# Row 42: update B42, C42, D42, E42
Set-TextValue "B42" "PaxDollar"
Set-TextValue "C42" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  -0.35%  "

# Row 43: update E43
Set-TextValue "E43" "  +0.11%  "

# Row 44: update D44, E44
Set-TextValue "D44" "0.764"
Set-TextValue "E44" "  +0.08%  "

# Row 45: update D45, E45
Set-TextValue "D45" "62.01"
Set-TextValue "E45" "  -0.70%  "

# Row 46: update D46, E46
Set-TextValue "D46" "1.728.13"
Set-TextValue "E46" "  +0.75%  "

# Row 47: update D47, E47
Set-TextValue "D47" "85.34"
Set-TextValue "E47" "  -4.08%  "

# Row 48: update E48
Set-TextValue "E48" "  -3.90%  "

# Row 49: update E49
Set-TextValue "E49" "  -0.60%  "

# Row 50: update D50, E50
Set-TextValue "D50" "0.0975"
Set-TextValue "E50" "  -2.80%  "

# Row 51: update D51, E51
Set-TextValue "D51" "0.998"
Set-TextValue "E51" "  -0.55%  "
